$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

# --- Update the running "counter" numbers shown inside the ellipse shapes ---
# First map() pass: ellipses with original numbers 3,4,5 become 2,2,3
$s.Shapes.Item(17).TextFrame.TextRange.Text = "2"   # Ellipse 24 : 3 -> 2
$s.Shapes.Item(18).TextFrame.TextRange.Text = "2"   # Ellipse 25 : 4 -> 2
$s.Shapes.Item(19).TextFrame.TextRange.Text = "3"   # Ellipse 26 : 5 -> 3

# Second aggregateByKey() pass: ellipses with original numbers 3,4,5,6,11 become 2,2,3,4,7
$s.Shapes.Item(22).TextFrame.TextRange.Text = "2"   # Ellipse 29 : 3 -> 2
$s.Shapes.Item(23).TextFrame.TextRange.Text = "2"   # Ellipse 30 : 4 -> 2
$s.Shapes.Item(24).TextFrame.TextRange.Text = "3"   # Ellipse 31 : 5 -> 3
$s.Shapes.Item(25).TextFrame.TextRange.Text = "4"   # Ellipse 32 : 6 -> 4
$s.Shapes.Item(26).TextFrame.TextRange.Text = "7"   # Ellipse 33 : 11 -> 7

# --- Add "Appear" click-entrance animations (with a build list) for all the counter ellipses ---
$seq = $s.TimeLine.MainSequence
for ($i = 15; $i -le 26; $i++) {
    $shp = $s.Shapes.Item($i)
    $eff = $seq.AddEffect($shp, 1, 0, 1)
}
